# Auto-generated from diff analysis: update scheduled-runner market data for Durandal Profits
$wb = $excel.ActiveWorkbook

# ALC row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 43479156
$ws.Range("I53").Value = 66667612
$ws.Range("J53").Value = 803.5
$ws.Range("K53").Value = 66667612
$ws.Range("L53").Value = 803.5
$ws.Range("M53").Value = -66666975
$ws.Range("N53").Value = -2077.5

# ALC row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6268.1816
$ws.Range("I141").Value = 6595
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 19785
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -14605
$ws.Range("N141").Value = -19360

# ARM row 37: Get Shirty / Steel Chainmail
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8200.333000000001
$ws.Range("J37").Value = 8200.333000000001
$ws.Range("L37").Value = 8200.333000000001
$ws.Range("N37").Value = -8746.333000000001

# ARM row 55: Employee Retention / Mythril Elmo
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 26200.092
$ws.Range("J55").Value = 26200.092
$ws.Range("L55").Value = 26200.092
$ws.Range("N55").Value = -26830.092

# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1671.4286
$ws.Range("I61").Value = 1671.4286
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1671.4286
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1459.4286
$ws.Range("N61").ClearContents()

# ARM row 80: A Squire to Inspire / Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 13099.75
$ws.Range("J80").Value = 16666.334
$ws.Range("L80").Value = 16666.334
$ws.Range("N80").Value = -18662.334

# ARM row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 13099.75
$ws.Range("J83").Value = 16666.334
$ws.Range("L83").Value = 49999.00199999999
$ws.Range("N83").Value = -59983.00199999999

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2170.3125
$ws.Range("I132").Value = 1309.9524
$ws.Range("J132").Value = 3812.818
$ws.Range("K132").Value = 3929.857199999999
$ws.Range("L132").Value = 11438.454
$ws.Range("M132").Value = -1399.857199999999
$ws.Range("N132").Value = -16498.454

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1671.4286
$ws.Range("I136").Value = 1671.4286
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5014.2858
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2464.2858
$ws.Range("N136").ClearContents()

# BSM row 20: Smelt and Dealt / Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3601.85
$ws.Range("I20").Value = 3456.4443
$ws.Range("J20").Value = 3720.818
$ws.Range("K20").Value = 3456.4443
$ws.Range("L20").Value = 3720.818
$ws.Range("M20").Value = -3209.4443
$ws.Range("N20").Value = -4214.818

# BSM row 82: Spirituality Inspector / Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 11823.909
$ws.Range("I82").Value = 3985.7144
$ws.Range("J82").Value = 25540.75
$ws.Range("K82").Value = 3985.7144
$ws.Range("L82").Value = 25540.75
$ws.Range("M82").Value = -3602.7144
$ws.Range("N82").Value = -26306.75

# BSM row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 11823.909
$ws.Range("I85").Value = 3985.7144
$ws.Range("J85").Value = 25540.75
$ws.Range("K85").Value = 3985.7144
$ws.Range("L85").Value = 25540.75
$ws.Range("M85").Value = -2659.7144
$ws.Range("N85").Value = -28192.75

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4667.1113
$ws.Range("I134").Value = 3429.1428
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 10287.4284
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -7752.428400000001
$ws.Range("N134").Value = -32070

# CRP row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9201.333000000001
$ws.Range("J50").Value = 9201.333000000001
$ws.Range("L50").Value = 9201.333000000001
$ws.Range("N50").Value = -10451.333

# CRP row 51: Greenstone for Greenhorns / Jade Crook
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9380
$ws.Range("J51").Value = 9380
$ws.Range("L51").Value = 9380
$ws.Range("N51").Value = -10852

# CRP row 59: Bow Down to Magic / Crab Bow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 13075.333
$ws.Range("J59").Value = 13075.333
$ws.Range("L59").Value = 13075.333
$ws.Range("N59").Value = -15365.333

# CRP row 60: Bowing to Greater Power / Yew Longbow
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 8200.666999999999
$ws.Range("J60").Value = 8200.666999999999
$ws.Range("L60").Value = 8200.666999999999
$ws.Range("N60").Value = -9222.666999999999

# CRP row 61: Incant Now, Think Later / Jade Crook
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9380
$ws.Range("J61").Value = 9380
$ws.Range("L61").Value = 9380
$ws.Range("N61").Value = -10076

# CRP row 109: Playing the Market / White Oak Necklace
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 16996.666
$ws.Range("J109").Value = 17995
$ws.Range("L109").Value = 17995
$ws.Range("N109").Value = -20075

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3935.5715
$ws.Range("I132").Value = 3058.8572
$ws.Range("J132").Value = 4812.2856
$ws.Range("K132").Value = 9176.571599999999
$ws.Range("L132").Value = 14436.8568
$ws.Range("M132").Value = -6646.571599999999
$ws.Range("N132").Value = -19496.8568

# CRP row 134: Wood You Be Quiet / Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2357.182
$ws.Range("I134").Value = 2014.3334
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 6043.0002
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -3508.0002
$ws.Range("N134").Value = -16770

# CUL row 113: Can't Eat Just One / Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 796.59
$ws.Range("I113").Value = 471.2
$ws.Range("J113").Value = 905.05334
$ws.Range("K113").Value = 1413.6
$ws.Range("L113").Value = 2715.16002
$ws.Range("M113").Value = 756.4000000000001
$ws.Range("N113").Value = -7055.16002

# CUL row 123: Topping Up the Pot / Zurek
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2187.125
$ws.Range("I123").Value = 2071
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 6213
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = -3763
$ws.Range("N123").Value = -13900

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7463547.5
$ws.Range("I131").Value = 900.9091
$ws.Range("J131").Value = 8929424
$ws.Range("K131").Value = 2702.7273
$ws.Range("L131").Value = 26788272
$ws.Range("M131").Value = 2337.2727
$ws.Range("N131").Value = -26798352

# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3251.389
$ws.Range("I137").Value = 1515
$ws.Range("J137").Value = 3468.4375
$ws.Range("K137").Value = 4545
$ws.Range("L137").Value = 10405.3125
$ws.Range("M137").Value = 555
$ws.Range("N137").Value = -20605.3125

# CUL row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 6581.4546
$ws.Range("I138").Value = 7912
$ws.Range("J138").Value = 3033.3333
$ws.Range("K138").Value = 23736
$ws.Range("L138").Value = 9099.999899999999
$ws.Range("M138").Value = -18596
$ws.Range("N138").Value = -19379.9999

# CUL row 139: Najoothie / Wild Banana Blend
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3992.1853
$ws.Range("I139").Value = 6223.875
$ws.Range("J139").Value = 3052.5264
$ws.Range("K139").Value = 18671.625
$ws.Range("L139").Value = 9157.5792
$ws.Range("M139").Value = -13531.625
$ws.Range("N139").Value = -19437.5792

# CUL row 141: Ocean Explosion / Acqua Pazza
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 9232.177
$ws.Range("I141").Value = 8368.375
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 25105.125
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -19925.125
$ws.Range("N141").Value = -40360

# GSM row 102: Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2387.4546
$ws.Range("I102").Value = 2258
$ws.Range("J102").Value = 2614
$ws.Range("K102").Value = 2258
$ws.Range("L102").Value = 2614
$ws.Range("M102").Value = -636
$ws.Range("N102").Value = -5858

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2834.1155
$ws.Range("I122").Value = 3193.0625
$ws.Range("J122").Value = 2259.8
$ws.Range("K122").Value = 9579.1875
$ws.Range("L122").Value = 6779.400000000001
$ws.Range("M122").Value = -7129.1875
$ws.Range("N122").Value = -11679.4

# GSM row 123: Workplace Workout / Ametrine Ring of Fending
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 17518.354
$ws.Range("J123").Value = 17518.354
$ws.Range("L123").Value = 17518.354
$ws.Range("N123").Value = -22418.354

# GSM row 132: On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4460.256
$ws.Range("I132").Value = 5016.0967
$ws.Range("J132").Value = 3024.3333
$ws.Range("K132").Value = 15048.2901
$ws.Range("L132").Value = 9072.999899999999
$ws.Range("M132").Value = -12518.2901
$ws.Range("N132").Value = -14132.9999
